# Atualizado por script em 05-11-2023 08:45
#
# 1) Swap the betting-odds data for the "AVS vs Vilaverdense" and
#    "FC Porto B vs Maritimo" matches (rows 39/40 keep the same Indice and
#    match date in columns A:E, only columns F:V are exchanged between the
#    two rows).
# 2) Append five new match rows (74-78), copying the formatting of the last
#    existing data row (73) for the styled columns (A = bold/centered/
#    bordered index, E = custom date/time format) so the new rows look like
#    the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MatchRow {
    param(
        [int]$Row,
        [string]$Home,
        [int]$HomeGoals,
        [string]$Away,
        [int]$AwayGoals,
        [double]$HomeOpenOdds,
        [string]$HomeOpenTime,
        [double]$HomeCloseOdds,
        [string]$HomeCloseTime,
        [double]$DrawOpenOdds,
        [string]$DrawOpenTime,
        [double]$DrawCloseOdds,
        [string]$DrawCloseTime,
        [double]$AwayOpenOdds,
        [string]$AwayOpenTime,
        [double]$AwayCloseOdds,
        [string]$AwayCloseTime,
        [string]$Url
    )

    $ws.Cells.Item($Row, 6).Value = $Home
    $ws.Cells.Item($Row, 7).Value = $HomeGoals
    $ws.Cells.Item($Row, 8).Value = $Away
    $ws.Cells.Item($Row, 9).Value = $AwayGoals

    $ws.Cells.Item($Row, 10).Value = $HomeOpenOdds
    $ws.Cells.Item($Row, 11).Value = $HomeOpenTime
    $ws.Cells.Item($Row, 12).Value = $HomeCloseOdds
    $ws.Cells.Item($Row, 13).Value = $HomeCloseTime

    $ws.Cells.Item($Row, 14).Value = $DrawOpenOdds
    $ws.Cells.Item($Row, 15).Value = $DrawOpenTime
    $ws.Cells.Item($Row, 16).Value = $DrawCloseOdds
    $ws.Cells.Item($Row, 17).Value = $DrawCloseTime

    $ws.Cells.Item($Row, 18).Value = $AwayOpenOdds
    $ws.Cells.Item($Row, 19).Value = $AwayOpenTime
    $ws.Cells.Item($Row, 20).Value = $AwayCloseOdds
    $ws.Cells.Item($Row, 21).Value = $AwayCloseTime

    $ws.Cells.Item($Row, 22).Value = $Url
}

function Set-NewRow {
    param(
        [int]$Row,
        [int]$Indice,
        [double]$MatchDate,
        [string]$Home,
        [int]$HomeGoals,
        [string]$Away,
        [int]$AwayGoals,
        [double]$HomeOpenOdds,
        [string]$HomeOpenTime,
        [double]$HomeCloseOdds,
        [string]$HomeCloseTime,
        [double]$DrawOpenOdds,
        [string]$DrawOpenTime,
        [double]$DrawCloseOdds,
        [string]$DrawCloseTime,
        [double]$AwayOpenOdds,
        [string]$AwayOpenTime,
        [double]$AwayCloseOdds,
        [string]$AwayCloseTime,
        [string]$Url
    )

    # Column A: bold/centered/bordered index - clone formatting from the
    # existing last row (73) so the new rows match the table exactly, then
    # overwrite the value.
    $ws.Range("A73").Copy()
    $ws.Range("A" + $Row).PasteSpecial(-4122)
    $ws.Cells.Item($Row, 1).Value = $Indice

    $ws.Cells.Item($Row, 2).Value = "portugal"
    $ws.Cells.Item($Row, 3).Value = "liga-portugal-2"
    $ws.Cells.Item($Row, 4).Value = "2023-2024"

    # Column E: match date/time - clone the custom date/time number format
    # from row 73, then overwrite the value.
    $ws.Range("E73").Copy()
    $ws.Range("E" + $Row).PasteSpecial(-4122)
    $ws.Cells.Item($Row, 5).Value = $MatchDate

    Set-MatchRow $Row $Home $HomeGoals $Away $AwayGoals `
        $HomeOpenOdds $HomeOpenTime $HomeCloseOdds $HomeCloseTime `
        $DrawOpenOdds $DrawOpenTime $DrawCloseOdds $DrawCloseTime `
        $AwayOpenOdds $AwayOpenTime $AwayCloseOdds $AwayCloseTime `
        $Url
}

# --- 1) Swap rows 39 and 40 (F:V only) -------------------------------------

Set-MatchRow 39 "FC Porto B" 0 "Maritimo" 2 `
    2.76 "13/09/2023 21:12" 2.49 "16/09/2023 11:00" `
    3.22 "13/09/2023 21:12" 3.34 "16/09/2023 10:30" `
    2.58 "13/09/2023 21:12" 3.01 "16/09/2023 11:00" `
    "https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-maritimo/tORXnMP1/"

Set-MatchRow 40 "AVS" 2 "Vilaverdense" 0 `
    1.75 "12/09/2023 12:12" 1.82 "16/09/2023 11:53" `
    3.71 "12/09/2023 12:12" 3.69 "16/09/2023 11:53" `
    4.89 "12/09/2023 12:12" 4.66 "16/09/2023 11:53" `
    "https://www.betexplorer.com/football/portugal/liga-portugal-2/avs-vilaverdense-fc/rVCSkOuq/"

# --- 2) Append new rows 74-78 -----------------------------------------------

Set-NewRow 74 73 45234.5 "FC Porto B" 2 "Feirense" 0 `
    1.98 "01/11/2023 16:12" 1.84 "04/11/2023 11:59" `
    3.57 "01/11/2023 16:12" 3.78 "04/11/2023 11:59" `
    3.87 "01/11/2023 16:12" 4.41 "04/11/2023 11:58" `
    "https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-feirense/jTL6QSDN/"

Set-NewRow 75 74 45234.5 "Mafra" 0 "Leixoes" 1 `
    1.88 "01/11/2023 16:12" 1.93 "04/11/2023 11:48" `
    3.73 "01/11/2023 16:12" 3.55 "04/11/2023 11:51" `
    3.8 "01/11/2023 16:12" 4.21 "04/11/2023 11:51" `
    "https://www.betexplorer.com/football/portugal/liga-portugal-2/mafra-leixoes/YiBBPnTT/"

Set-NewRow 76 75 45234.6875 "Os Belenenses" 0 "Vilaverdense" 2 `
    2.31 "01/11/2023 16:12" 2.24 "04/11/2023 16:21" `
    3.34 "01/11/2023 16:12" 3.34 "04/11/2023 16:21" `
    3.05 "01/11/2023 16:12" 3.47 "04/11/2023 16:21" `
    "https://www.betexplorer.com/football/portugal/liga-portugal-2/cf-os-belenenses-vilaverdense-fc/r3IbSlqB/"

Set-NewRow 77 76 45234.79166666666 "Benfica B" 1 "Penafiel" 0 `
    1.91 "29/10/2023 16:42" 2.26 "04/11/2023 18:53" `
    3.67 "29/10/2023 16:42" 3.59 "04/11/2023 18:53" `
    4.01 "29/10/2023 16:42" 3.2 "04/11/2023 18:53" `
    "https://www.betexplorer.com/football/portugal/liga-portugal-2/benfica-penafiel/xjmbUAEb/"

Set-NewRow 78 77 45234.79166666666 "Nacional" 1 "Santa Clara" 1 `
    2.98 "01/11/2023 16:12" 2.81 "04/11/2023 18:58" `
    3.27 "01/11/2023 16:12" 3.23 "04/11/2023 18:52" `
    2.39 "01/11/2023 16:12" 2.72 "04/11/2023 18:52" `
    "https://www.betexplorer.com/football/portugal/liga-portugal-2/nacional-santa-clara/xQH2R8bH/"

Write-Output "Applied odds swap for rows 39-40 and appended rows 74-78"
